$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, bordered, centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New column values, per row
$values = @{
    2  = @(8, 8)
    3  = @(9, 9)
    4  = @(7, 7)
    5  = @(9, 9)
    6  = @(9, 9)
    7  = @(6, 7)
    8  = @(9, 9)
    9  = @(5, 5)
    10 = @(6, 6)
    11 = @(7, 7)
    12 = @(7, 7)
    13 = @(5, 6)
    14 = @(5, 5)
    15 = @(5, 5)
    16 = @(8, 8)
    17 = @(8, 8)
    18 = @(8, 8)
    19 = @(7, 7)
    20 = @(5, 5)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

Write-Output "done"
